# Updated cryptos list on Mon Feb  5 20:00:33 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row on the
# active sheet, plus the HuobiToken / RocketPoolETH row swap (rows 48-49).
#
# D-column values that look like plain decimals (single '.') are written
# with a leading apostrophe so Excel stores them as text (preserving
# trailing zeros / exact formatting) instead of silently converting them to
# numbers. Values that already contain two '.' (thousands separator style,
# e.g. "42.804.05") are unambiguous as text and are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.804.05"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "2.305.60"
$ws.Range("E3").Value = "  +1.09%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'301.62"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").Value = "'96.08"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").Value = "'0.510"
$ws.Range("E7").Value = "  +0.58%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("D10").Value = "'34.61"
$ws.Range("E10").Value = "  -2.46%  "

$ws.Range("D11").Value = "'19.18"
$ws.Range("E11").Value = "  +4.88%  "

$ws.Range("D12").Value = "'0.0789"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("D14").Value = "'6.79"
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").Value = "2.667.01"
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("D16").Value = "2.309.88"
$ws.Range("E16").Value = "  +1.12%  "

$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "42.742.67"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("E19").Value = "  -5.73%  "

$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").Value = "'67.78"

$ws.Range("E23").Value = "  +6.67%  "

$ws.Range("D24").Value = "'235.50"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = "  -1.49%  "

$ws.Range("D27").Value = "'24.39"
$ws.Range("E27").Value = "  -2.92%  "

$ws.Range("E28").Value = "  +14.55%  "

$ws.Range("D29").Value = "'165.62"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("D30").Value = "'9.07"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("D31").Value = "'32.06"
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("D34").Value = "'17.60"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E35").Value = "  -6.44%  "

$ws.Range("D36").Value = "'0.0700"
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("E37").Value = "  -2.70%  "

$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("D42").Value = "'20.24"
$ws.Range("E42").Value = "  +11.48%  "

$ws.Range("D43").Value = "1.967.60"
$ws.Range("E43").Value = "  -1.62%  "

$ws.Range("E44").Value = "  +5.18%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("E47").Value = "  -0.32%  "

# Row 48/49: HuobiToken and RocketPoolETH swap places.
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.530.86"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.83"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").Value = "'53.25"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").Value = "'71.46"
$ws.Range("E51").Value = "  +0.21%  "
